# Libra vendor-upload template: add an "Is IC?" column (is_ic) to the
# "vendor" sheet, after the existing "ASCII Vendor Code" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendor")
$ws.Activate() | Out-Null

# New column L: machine key in row 1, human label in row 2 (matches the
# existing header/label pattern used by every other column on this sheet).
$ws.Range("L1").Value = "is_ic"
$ws.Range("L2").Value = "Is IC?"

# Re-fit the columns whose best-fit width shifts now that row 1 holds a
# wider mix of header keys (A:D were impacted by the new, wider content).
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null

# Leave the same kind of selection the author ended up with after editing.
$ws.Range("G2").Select() | Out-Null
